$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure D and E columns (price / volume change) stay formatted as text,
# matching the original inlineStr cell types, so values such as "225.70"
# or "1.004" are not coerced into numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '27.560.45'
$ws.Range('E2').Value = '  +5.47%  '

$ws.Range('D3').Value = '1.726.26'
$ws.Range('E3').Value = '  +4.41%  '

$ws.Range('E4').Value = '  +0.07%  '

$ws.Range('D5').Value = '225.70'
$ws.Range('E5').Value = '  +3.22%  '

$ws.Range('D6').Value = '0.5378'
$ws.Range('E6').Value = '  +3.02%  '

$ws.Range('D7').Value = '1.004'
$ws.Range('E7').Value = '  +0.03%  '

$ws.Range('D8').Value = '0.2676'
$ws.Range('E8').Value = '  +0.89%  '

$ws.Range('D9').Value = '0.06614'
$ws.Range('E9').Value = '  +4.10%  '

$ws.Range('E10').Value = '  +6.75%  '

$ws.Range('D11').Value = '0.07713'
$ws.Range('E11').Value = '  +0.30%  '

$ws.Range('D12').Value = '4.613'
$ws.Range('E12').Value = '  -0.49%  '

$ws.Range('D13').Value = '1.723.36'
$ws.Range('E13').Value = '  +5.83%  '

$ws.Range('D14').Value = '1.964.30'
$ws.Range('E14').Value = '  +4.38%  '

$ws.Range('D15').Value = '0.5855'
$ws.Range('E15').Value = '  +4.44%  '

$ws.Range('D16').Value = '0.0₅8324'
$ws.Range('E16').Value = '  +1.59%  '

$ws.Range('D17').Value = '68.05'
$ws.Range('E17').Value = '  +3.80%  '

$ws.Range('D18').Value = '27.570.56'
$ws.Range('E18').Value = '  +5.52%  '

$ws.Range('D19').Value = '221.99'
$ws.Range('E19').Value = '  +15.55%  '

$ws.Range('E20').Value = '  +0.08%  '

$ws.Range('D21').Value = '4.742'
$ws.Range('E21').Value = '  +2.08%  '

$ws.Range('D22').Value = '10.68'
$ws.Range('E22').Value = '  +1.82%  '

$ws.Range('E23').Value = '  +2.71%  '

$ws.Range('E24').Value = '  +0.07%  '

$ws.Range('D25').Value = '148.44'
$ws.Range('E25').Value = '  +1.98%  '

$ws.Range('B26').Value = 'Toncoin'
$ws.Range('C26').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D26').Value = '1.704'
$ws.Range('E26').Value = '  +12.65%  '

$ws.Range('B27').Value = 'Stellar'
$ws.Range('C27').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D27').Value = '0.1238'
$ws.Range('E27').Value = '  +3.62%  '

$ws.Range('D28').Value = '7.415'
$ws.Range('E28').Value = '  +2.35%  '

$ws.Range('D29').Value = '16.68'
$ws.Range('E29').Value = '  +4.76%  '

$ws.Range('D30').Value = '0.05585'
$ws.Range('E30').Value = '  +1.75%  '

$ws.Range('D31').Value = '1.303'
$ws.Range('E31').Value = '  +2.53%  '

$ws.Range('D32').Value = '3.553'
$ws.Range('E32').Value = '  +2.74%  '

$ws.Range('D33').Value = '3.457'
$ws.Range('E33').Value = '  +2.36%  '

$ws.Range('E34').Value = '  +6.59%  '

$ws.Range('D35').Value = '0.9655'
$ws.Range('E35').Value = '  +1.22%  '

$ws.Range('D36').Value = '2.825'
$ws.Range('E36').Value = '  +1.39%  '

$ws.Range('E37').Value = '  +2.17%  '

$ws.Range('D38').Value = '0.5960'
$ws.Range('E38').Value = '  +5.58%  '

$ws.Range('D39').Value = '0.01651'
$ws.Range('E39').Value = '  +4.34%  '

$ws.Range('D40').Value = '5.931'
$ws.Range('E40').Value = '  +1.06%  '

$ws.Range('D41').Value = '0.8579'
$ws.Range('E41').Value = '  +2.69%  '

$ws.Range('D42').Value = '1.056.64'
$ws.Range('E42').Value = '  +2.80%  '

$ws.Range('E43').Value = '  +0.08%  '

$ws.Range('D44').Value = '101.52'
$ws.Range('E44').Value = '  +0.25%  '

$ws.Range('D45').Value = '1.869.26'
$ws.Range('E45').Value = '  +4.25%  '

$ws.Range('D46').Value = '0.0₈114'
$ws.Range('E46').Value = '  +6.01%  '

$ws.Range('D47').Value = '59.24'
$ws.Range('E47').Value = '  +2.55%  '

$ws.Range('D48').Value = '8.211'
$ws.Range('E48').Value = '  +2.88%  '

$ws.Range('D49').Value = '0.4440'
$ws.Range('E49').Value = '  +2.31%  '

$ws.Range('D50').Value = '1.001'
$ws.Range('E50').Value = '  +0.09%  '

$ws.Range('D51').Value = '0.05268'
$ws.Range('E51').Value = '  +1.46%  '

# Restore default (unstyled) appearance now that values are set, so no
# stray number-format style gets attached to the cells.
$ws.Range("D2:E51").Style = "Normal"

Write-Output "Updated cryptos list"
